$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the template row (row 4) down into the two new
# rows (5 and 6) so the new data picks up the same cell styles used by
# the other data rows in this report. Row 5 only gets columns A:J (it
# never had a K cell), row 6 reuses the existing A:K template (its K
# cell stays present, but empty, same as before).
$ws.Range("A4:J4").Copy() | Out-Null
$ws.Range("A5:J5").PasteSpecial(-4122) | Out-Null
$ws.Range("A4:K4").Copy() | Out-Null
$ws.Range("A6:K6").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Row 5 - "Crumpet" facility
$ws.Range("A5").Value = "Crumpet GEF"
$ws.Range("B5").Value = 20001371
$ws.Range("C5").Value = "Crumpet exporter"
$ws.Range("D5").Value = "GBP"
$ws.Range("E5").Value = 7000000
$ws.Range("F5").Value = 3938753.8
$ws.Range("G5").Value = 777
$ws.Range("H5").Value = 456
$ws.Range("I5").Value = "GBP"
$ws.Range("J5").Value = "GBP"

# Row 6 - "Scone" facility
$ws.Range("A6").Value = "Scone GEF"
$ws.Range("B6").Value = 20001371
$ws.Range("C6").Value = "Scone exporter"
$ws.Range("D6").Value = "GBP"
$ws.Range("E6").Value = 770000
$ws.Range("F6").Value = 761579.37
$ws.Range("G6").Value = 777
$ws.Range("H6").Value = 456.77
$ws.Range("I6").Value = "GBP"
$ws.Range("J6").Value = "GBP"
$ws.Range("K6").ClearContents()
